# Update "想去人数" (interested-people count) figures that changed between
# the two scraped snapshots of the 南宁-漫展信息 workbook.
#
# Sheet "展览" (Exhibitions) - rows keyed by F column values:
#   F2: 14813 -> 14953
#   F4: 704   -> 705
#   F5: 245   -> 246
#   F6: 618   -> 624
#   F7: 1578  -> 1592
#   F8: 147   -> 148
#
# Sheet "全部类型" (All types) - same events, different row numbers:
#   F2:  14813 -> 14954
#   F4:  704   -> 705
#   F5:  245   -> 246
#   F8:  618   -> 624
#   F9:  1578  -> 1592
#   F11: 147   -> 148

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 14953
$wsExhibition.Range("F4").Value = 705
$wsExhibition.Range("F5").Value = 246
$wsExhibition.Range("F6").Value = 624
$wsExhibition.Range("F7").Value = 1592
$wsExhibition.Range("F8").Value = 148

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F2").Value = 14954
$wsAllTypes.Range("F4").Value = 705
$wsAllTypes.Range("F5").Value = 246
$wsAllTypes.Range("F8").Value = 624
$wsAllTypes.Range("F9").Value = 1592
$wsAllTypes.Range("F11").Value = 148
